# Test data added for Italy
$wb = $excel.ActiveWorkbook

# Duplicate the "Slovakia" sheet (our template) and place the copy after the
# last sheet in the workbook; this gives the new sheet the same layout,
# column widths, merged cells and styles as the existing country sheets.
$template = $wb.Worksheets.Item("Slovakia")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy([System.Reflection.Missing]::Value, $lastSheet)

$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "Italy"

# Fill in the Italy-specific values.
$ws.Range("B2").Value = "Italy Market"
$ws.Range("B4").Value = "NGC-3145/T2155"

# The ticket cell picked up a slightly larger, borderless font when it was
# authored - match that formatting.
$ws.Range("B4").Font.Size = 12
$ws.Range("B4").Font.Color = 0
$ws.Range("B4").Borders.LineStyle = -4142
$ws.Rows.Item(4).RowHeight = 15.6

# The new Italy sheet is the one the user finished on.
$ws.Activate()
